# Penalty/Reward System rework (unfinished): remove some weekly/monthly
# rows from the PO data so the remaining rows shift up.

$wb = $excel.ActiveWorkbook

# --- "Weekly Quantity" sheet: delete rows 4-6 (2023-07-09, 2023-07-16, 2023-07-23) ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws1.Rows("4:6").Delete()

# --- "Monthly Trend" sheet: delete row 3 (2023-08-01) ---
$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Rows("3:3").Delete()
